$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: "Property" -> "PostName"
$ws.Range("A1").Value = "PostName"

# B1: new -> "BlogName"
$ws.Range("B1").Value = "BlogName"

# A2 stays "PropsGehenRaus" (unchanged content, but shared string index shifts automatically)

# B2: empty -> "KenBlock" (keep existing style s="1")
$ws.Range("B2").Value = "KenBlock"

# Move selection to B3
$ws.Range("B3").Select()
